$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.7597
$ws.Range("B7").Value = 4.717699999999997
$ws.Range("A8").Value = -22.28100000000001
$ws.Range("A10").Value = -21.8271
$ws.Range("A12").Value = -21.5636
$ws.Range("B15").Value = 5.059099999999995
$ws.Range("A18").Value = -21.87279999999999
$ws.Range("B18").Value = 6.604000000000001
$ws.Range("C18").Value = -11.8307
$ws.Range("C19").Value = -11.3514
$ws.Range("B20").Value = 8.864400000000003
$ws.Range("C27").Value = -13.1888
$ws.Range("B29").Value = 4.937600000000003
$ws.Range("B30").Value = 5.050100000000001
$ws.Range("B31").Value = 5.0898
$ws.Range("C31").Value = -13.3054
$ws.Range("A37").Value = -19.49369999999999
$ws.Range("C38").Value = -13.1658
$ws.Range("B40").Value = 9.144399999999989
$ws.Range("C42").Value = -12.12560000000001
$ws.Range("C44").Value = -13.66309999999998
$ws.Range("C47").Value = -12.2308
$ws.Range("B50").Value = 4.611000000000001
$ws.Range("A55").Value = -22.4148
$ws.Range("C58").Value = -12.3502
$ws.Range("C65").Value = -12.1855
$ws.Range("A68").Value = -21.45170000000001
$ws.Range("B68").Value = 4.491999999999998
$ws.Range("C73").Value = -12.2189
$ws.Range("B76").Value = 6.188000000000001
$ws.Range("A77").Value = -19.97609999999998
$ws.Range("A78").Value = -20.15919999999998
$ws.Range("A81").Value = -21.9575
$ws.Range("A82").Value = -21.91449999999999
$ws.Range("B87").Value = 5.241599999999994
$ws.Range("B88").Value = 4.391299999999998
$ws.Range("C90").Value = -13.0379
$ws.Range("C94").Value = -10.10750000000001
$ws.Range("C95").Value = -11.8181
$ws.Range("B96").Value = 5.131400000000008
$ws.Range("B98").Value = 6.7342
$ws.Range("B101").Value = 9.257999999999988
$ws.Range("C101").Value = -13.42320000000001
$ws.Range("B102").Value = 8.600700000000005
